# Add a new "2021年" data row (row 5) to Sheet1, mirroring the existing
# 2018-2020 rows. Values taken from the underlying published dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = '2021年'

$ws.Range("B5").Value = 14.8
$ws.Range("C5").Value = 38.3
$ws.Range("E5").Value = 19.2
$ws.Range("F5").Value = 31.1
$ws.Range("H5").Value = -3.5
$ws.Range("I5").Value = 20.9
$ws.Range("J5").Value = 7.4
$ws.Range("K5").Value = 11.1
$ws.Range("M5").Value = 20.7
$ws.Range("O5").Value = -4.6
$ws.Range("P5").Value = 11.3
$ws.Range("Q5").Value = -59.1
$ws.Range("S5").Value = 117.1
$ws.Range("T5").Value = 30.7
$ws.Range("U5").Value = 18.2
$ws.Range("V5").Value = -19.6
$ws.Range("W5").Value = 35.2
$ws.Range("X5").Value = 19.5
$ws.Range("Y5").Value = 23.6
$ws.Range("Z5").Value = -35.9
$ws.Range("AA5").Value = 18.5
$ws.Range("AB5").Value = -43.1
$ws.Range("AC5").Value = -50.7
$ws.Range("AD5").Value = -29.7
$ws.Range("AE5").Value = 92.7
$ws.Range("AK5").Value = 2.4
$ws.Range("AL5").Value = -31.3
$ws.Range("AM5").Value = 51.5
$ws.Range("AN5").Value = -21.0
$ws.Range("AO5").Value = 5.7
$ws.Range("AP5").Value = -68.5
$ws.Range("AQ5").Value = 12.7
$ws.Range("AR5").Value = -49.5
$ws.Range("AU5").Value = -87.7
$ws.Range("AV5").Value = 10.7
$ws.Range("AW5").Value = -49.5
$ws.Range("AX5").Value = 96.9
$ws.Range("AY5").Value = 16.7
$ws.Range("AZ5").Value = 71.7
$ws.Range("BA5").Value = -14.9
$ws.Range("BB5").Value = 146.4
$ws.Range("BC5").Value = 42.5
$ws.Range("BE5").Value = -1.5
$ws.Range("BF5").Value = 36.1
$ws.Range("BG5").Value = 7.9
$ws.Range("BH5").Value = 93.6
$ws.Range("BI5").Value = -44.1
$ws.Range("BJ5").Value = 10.4
$ws.Range("BK5").Value = 64.8
$ws.Range("BL5").Value = -6.8
$ws.Range("BN5").Value = 5.4
$ws.Range("BO5").Value = 0.1
$ws.Range("BP5").Value = 221.7
$ws.Range("BQ5").Value = -71.5
$ws.Range("BR5").Value = 20.1
$ws.Range("BS5").Value = 25.7
$ws.Range("BT5").Value = -29.4
$ws.Range("BU5").Value = 62.4
$ws.Range("BV5").Value = -9.2
$ws.Range("BW5").Value = -16.7
$ws.Range("BX5").Value = 38.3
$ws.Range("BY5").Value = 37.2
$ws.Range("BZ5").Value = 22.3
$ws.Range("CA5").Value = 115.1
$ws.Range("CB5").Value = 125.3
$ws.Range("CC5").Value = 21.9
$ws.Range("CE5").Value = -7.0
$ws.Range("CF5").Value = -23.5
$ws.Range("CG5").Value = -53.5
$ws.Range("CH5").Value = -71.7
$ws.Range("CI5").Value = 91.1
$ws.Range("CJ5").Value = 173.1
$ws.Range("CK5").Value = 4.5
$ws.Range("CL5").Value = -8.8
$ws.Range("CN5").Value = 15.2
$ws.Range("CO5").Value = 28.4
$ws.Range("CP5").Value = 9.5
$ws.Range("CQ5").Value = -34.1
$ws.Range("CR5").Value = -89.7
$ws.Range("CS5").Value = -2.0
$ws.Range("CT5").Value = 20.4
$ws.Range("CU5").Value = 36.7
$ws.Range("CV5").Value = 15.9
$ws.Range("CW5").Value = 232.2
$ws.Range("CX5").Value = -36.4
$ws.Range("CY5").Value = 112.9
$ws.Range("CZ5").Value = 57.5
$ws.Range("DA5").Value = 19.1
$ws.Range("DB5").Value = -33.2
$ws.Range("DC5").Value = 142.8
$ws.Range("DE5").Value = -9.9
$ws.Range("DF5").Value = 37.3
$ws.Range("DG5").Value = 67.3
$ws.Range("DH5").Value = 29.8
$ws.Range("DI5").Value = -54.1
$ws.Range("DJ5").Value = 8.4
$ws.Range("DK5").Value = -85.0

# Copy the bold/centered/bordered style already used for the year labels
# in A2:A4 onto the new label cell.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
